$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 2 (shifts existing rows 2..167 down to 3..168)
$ws.Rows.Item(2).Insert()

# 2. Copy formats from row 3 (the old row 2, now shifted down) into the new row 2
#    so it matches the rest of the data rows (center aligned, non-bold, 0.000 Basic Price format)
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Populate new row 2: same data as old row 2 (now row 3) except the Date column,
#    which advances to the new latest date (no new circular was published yet, so
#    Basic Price / Circular Date / Circular Link roll forward unchanged).
$ws.Range("A2").Value = "25-11-2025"
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 283
$ws.Range("E2").Value = "22.11.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-22-november-2025.pdf"

# 4. Rebuild the Hyperlinks collection: the engine does not shift hyperlink anchors
#    when rows are inserted, so clear them all and re-add one per non-empty Circular Link cell.
$ws.Cells.Hyperlinks.Delete()
for ($r = 2; $r -le 168; $r++) {
  $lnk = $ws.Cells.Item($r, 6).Value()
  if ($lnk) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $lnk)
  }
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()